$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new tenant data
$ws.Range("B2").Value = "Jibin Roy"
$ws.Range("C2").Value = "8111849588"
$ws.Range("D2").Value = "jibinroy949@gmail.com"
$ws.Range("E2").Value = "qwerty"
$ws.Range("F2").Value = "Aug-2024"
$ws.Range("G2").Value = "Single Private"
$ws.Range("H2").Value = "1"
$ws.Range("I2").Value = "15000"
$ws.Range("J2").Value = "Monthly Rent"
$ws.Range("K2").Value = "UPI"
$ws.Range("L2").Value = "Paid"
$ws.Range("M2").Value = "07-12-2024"

# Delete row 3 entirely (previously Dean Winchestor's record)
$ws.Rows.Item(3).Delete()

# Turn off right-to-left view
$ws.Activate()
$excel.ActiveWindow.DisplayRightToLeft = $false
